# Apply the "update database and change read_price algorithm" edit:
# All quarterly income-statement figures (rows 11-27, columns D:M) on the
# Overview sheet are reset to 0, except for the two rows whose figure is
# represented as a dash ("-") placeholder (rows 15 and 23).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

$dashRows = @(15, 23)

for ($row = 11; $row -le 27; $row++) {
    for ($col = 4; $col -le 13; $col++) {
        if ($dashRows -contains $row) {
            $ws.Cells.Item($row, $col).Value = "-"
        } else {
            $ws.Cells.Item($row, $col).Value = 0
        }
    }
}
